$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pedidos")

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 1) "P0001"
Set-TextValue $ws.Cells.Item(2, 2) "20498176144"
Set-TextValue $ws.Cells.Item(3, 1) "P0002"
Set-TextValue $ws.Cells.Item(3, 2) "10091942165"
Set-TextValue $ws.Cells.Item(4, 1) "P0003"
Set-TextValue $ws.Cells.Item(4, 2) "10458834445"
Set-TextValue $ws.Cells.Item(5, 1) "P0004"
Set-TextValue $ws.Cells.Item(5, 2) "10475747271"
Set-TextValue $ws.Cells.Item(6, 1) "P0005"
Set-TextValue $ws.Cells.Item(6, 2) "10090993998"
Set-TextValue $ws.Cells.Item(7, 1) "P0006"
Set-TextValue $ws.Cells.Item(7, 2) "10094512561"
Set-TextValue $ws.Cells.Item(8, 1) "P0007"
Set-TextValue $ws.Cells.Item(8, 2) "10096861848"
Set-TextValue $ws.Cells.Item(9, 1) "P0008"
Set-TextValue $ws.Cells.Item(9, 2) "10074990539"
Set-TextValue $ws.Cells.Item(10, 1) "P0009"
Set-TextValue $ws.Cells.Item(10, 2) "20508781475"
Set-TextValue $ws.Cells.Item(11, 1) "P0010"
Set-TextValue $ws.Cells.Item(11, 2) "20508781475"
Set-TextValue $ws.Cells.Item(12, 1) "P0011"
Set-TextValue $ws.Cells.Item(12, 2) "20508781475"
Set-TextValue $ws.Cells.Item(13, 1) "P0012"
Set-TextValue $ws.Cells.Item(13, 2) "20601560217"
Set-TextValue $ws.Cells.Item(14, 1) "P0013"
Set-TextValue $ws.Cells.Item(14, 2) "20545786959"
Set-TextValue $ws.Cells.Item(15, 1) "P0014"
Set-TextValue $ws.Cells.Item(15, 2) "10431424202"
Set-TextValue $ws.Cells.Item(16, 1) "P0015"
Set-TextValue $ws.Cells.Item(16, 2) "10459514991"
Set-TextValue $ws.Cells.Item(17, 1) "P0016"
Set-TextValue $ws.Cells.Item(17, 2) "20601560217"
Set-TextValue $ws.Cells.Item(18, 1) "P0017"
Set-TextValue $ws.Cells.Item(18, 2) "20509909157"
Set-TextValue $ws.Cells.Item(19, 1) "P0018"
Set-TextValue $ws.Cells.Item(19, 2) "10091942165"
Set-TextValue $ws.Cells.Item(20, 1) "P0019"
Set-TextValue $ws.Cells.Item(20, 2) "20600997433"
Set-TextValue $ws.Cells.Item(21, 1) "P0020"
Set-TextValue $ws.Cells.Item(21, 2) "20518051785"
Set-TextValue $ws.Cells.Item(22, 1) "P0021"
Set-TextValue $ws.Cells.Item(22, 2) "10073316206"
Set-TextValue $ws.Cells.Item(23, 1) "P0022"
Set-TextValue $ws.Cells.Item(23, 2) "20524524067"
Set-TextValue $ws.Cells.Item(24, 1) "P0023"
Set-TextValue $ws.Cells.Item(24, 2) "10086290796"
Set-TextValue $ws.Cells.Item(25, 1) "P0024"
Set-TextValue $ws.Cells.Item(25, 2) "20602510191"
Set-TextValue $ws.Cells.Item(26, 1) "P0025"
Set-TextValue $ws.Cells.Item(26, 2) "20516044820"
Set-TextValue $ws.Cells.Item(27, 1) "P0026"
Set-TextValue $ws.Cells.Item(27, 2) "20600093020"
Set-TextValue $ws.Cells.Item(28, 1) "P0027"
Set-TextValue $ws.Cells.Item(28, 2) "10104060400"
Set-TextValue $ws.Cells.Item(29, 1) "P0028"
Set-TextValue $ws.Cells.Item(29, 2) "20524524067"
Set-TextValue $ws.Cells.Item(30, 1) "P0029"
Set-TextValue $ws.Cells.Item(30, 2) "10091942165"
Set-TextValue $ws.Cells.Item(31, 1) "P0030"
Set-TextValue $ws.Cells.Item(31, 2) "10086290796"
Set-TextValue $ws.Cells.Item(32, 1) "P0031"
Set-TextValue $ws.Cells.Item(32, 2) "20548553386"
Set-TextValue $ws.Cells.Item(33, 1) "P0032"
Set-TextValue $ws.Cells.Item(33, 2) "10095279088"
Set-TextValue $ws.Cells.Item(34, 1) "P0033"
Set-TextValue $ws.Cells.Item(34, 2) "10429182153"
Set-TextValue $ws.Cells.Item(35, 1) "P0034"
Set-TextValue $ws.Cells.Item(35, 2) "20509909157"
Set-TextValue $ws.Cells.Item(36, 1) "P0035"
Set-TextValue $ws.Cells.Item(36, 2) "10091942165"
Set-TextValue $ws.Cells.Item(37, 1) "P0036"
Set-TextValue $ws.Cells.Item(37, 2) "10086290796"
Set-TextValue $ws.Cells.Item(38, 1) "P0037"
Set-TextValue $ws.Cells.Item(38, 2) "10404885541"
Set-TextValue $ws.Cells.Item(39, 1) "P0038"
Set-TextValue $ws.Cells.Item(39, 2) "20537112515"
Set-TextValue $ws.Cells.Item(40, 1) "P0039"
Set-TextValue $ws.Cells.Item(40, 2) "10453830816"
Set-TextValue $ws.Cells.Item(41, 1) "P0040"
Set-TextValue $ws.Cells.Item(41, 2) "10415543854"
Set-TextValue $ws.Cells.Item(42, 1) "P0041"
Set-TextValue $ws.Cells.Item(42, 2) "20478203676"
Set-TextValue $ws.Cells.Item(43, 1) "P0042"
Set-TextValue $ws.Cells.Item(43, 2) "10093662151"
Set-TextValue $ws.Cells.Item(44, 1) "P0043"
Set-TextValue $ws.Cells.Item(44, 2) "10091942165"
Set-TextValue $ws.Cells.Item(45, 1) "P0044"
Set-TextValue $ws.Cells.Item(45, 2) "20566486408"
Set-TextValue $ws.Cells.Item(46, 1) "P0045"
Set-TextValue $ws.Cells.Item(46, 2) "20602713831"
Set-TextValue $ws.Cells.Item(47, 1) "P0046"
Set-TextValue $ws.Cells.Item(47, 2) "10108095721"
Set-TextValue $ws.Cells.Item(48, 1) "P0047"
Set-TextValue $ws.Cells.Item(48, 2) "20566486408"
Set-TextValue $ws.Cells.Item(49, 1) "P0048"
Set-TextValue $ws.Cells.Item(49, 2) "10412628433"
Set-TextValue $ws.Cells.Item(50, 1) "P0049"
Set-TextValue $ws.Cells.Item(50, 2) "20546732305"
Set-TextValue $ws.Cells.Item(51, 1) "P0050"
Set-TextValue $ws.Cells.Item(51, 2) "10221033413"
Set-TextValue $ws.Cells.Item(52, 1) "P0051"
Set-TextValue $ws.Cells.Item(52, 2) "10200543675"
Set-TextValue $ws.Cells.Item(53, 1) "P0052"
Set-TextValue $ws.Cells.Item(53, 2) "10200543675"
Set-TextValue $ws.Cells.Item(54, 1) "P0053"
Set-TextValue $ws.Cells.Item(54, 2) "10486819192"
Set-TextValue $ws.Cells.Item(55, 1) "P0054"
Set-TextValue $ws.Cells.Item(55, 2) "20600093020"
Set-TextValue $ws.Cells.Item(56, 1) "P0055"
Set-TextValue $ws.Cells.Item(56, 2) "20601993091"
Set-TextValue $ws.Cells.Item(57, 1) "P0056"
Set-TextValue $ws.Cells.Item(57, 2) "10159480424"
Set-TextValue $ws.Cells.Item(58, 1) "P0057"
Set-TextValue $ws.Cells.Item(58, 2) "10152051404"
Set-TextValue $ws.Cells.Item(59, 1) "P0058"
Set-TextValue $ws.Cells.Item(59, 2) "10103382136"
Set-TextValue $ws.Cells.Item(60, 1) "P0059"
Set-TextValue $ws.Cells.Item(60, 2) "10100210598"
Set-TextValue $ws.Cells.Item(61, 1) "P0060"
Set-TextValue $ws.Cells.Item(61, 2) "10199657564"
Set-TextValue $ws.Cells.Item(62, 1) "P0061"
Set-TextValue $ws.Cells.Item(62, 2) "20498176144"
Set-TextValue $ws.Cells.Item(63, 1) "P0062"
Set-TextValue $ws.Cells.Item(63, 2) "20512002090"
Set-TextValue $ws.Cells.Item(64, 1) "P0063"
Set-TextValue $ws.Cells.Item(64, 2) "10715246011"
Set-TextValue $ws.Cells.Item(65, 1) "P0064"
Set-TextValue $ws.Cells.Item(65, 2) "17499553081"
Set-TextValue $ws.Cells.Item(66, 1) "P0065"
Set-TextValue $ws.Cells.Item(66, 2) "10152051404"
Set-TextValue $ws.Cells.Item(67, 1) "P0066"
Set-TextValue $ws.Cells.Item(67, 2) "10100210598"
Set-TextValue $ws.Cells.Item(68, 1) "P0067"
Set-TextValue $ws.Cells.Item(68, 2) "10096861848"
Set-TextValue $ws.Cells.Item(69, 1) "P0068"
Set-TextValue $ws.Cells.Item(69, 2) "20507349197"
Set-TextValue $ws.Cells.Item(70, 1) "P0069"
Set-TextValue $ws.Cells.Item(70, 2) "10100210598"
Set-TextValue $ws.Cells.Item(71, 1) "P0070"
Set-TextValue $ws.Cells.Item(71, 2) "10431424202"

# Row 45 column B gets an extra left-alignment style in the source edit
$ws.Cells.Item(45, 2).HorizontalAlignment = -4131

# Apply AutoFilter over the full data range and register the hidden _FilterDatabase name
$filterRange = $ws.Range("A1:H71")
$filterRange.AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", $filterRange)
$filterName.Visible = $false

# Restore the active selection to B1 as in the saved workbook
$ws.Range("B1").Select()
